# Automatische test-sync: 2025-07-27 19:12:50
# Appends a new log row (row 3) to the "Logs" sheet, mirroring row 2's
# mail-test content but with its own timestamp + follow-up flags, extends
# the conditional formatting ranges to cover the new row, and bumps the
# "Dashboard" sheet's tally for the "Overig" category.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- New row of data on "Logs" ---------------------------------------
$antwoord = "Geachte klant," + "`n" `
    + "Dank u voor uw bericht. Om u zo goed mogelijk van dienst te zijn, zou u wat meer details kunnen geven over wat u precies geregeld wilt hebben? Zo kunnen wij u beter helpen." + "`n" `
    + "Met vriendelijke groet," + "`n" `
    + "[Bedrijfsnaam] E-mailassistent"

$logs.Range("A3").Value = "Kun jij dit even regelen?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D3").Value = "Overig"
$logs.Range("E3").Value = $antwoord
$logs.Range("F3").Value = "2025-07-27 19:12:17"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Nee"
$logs.Range("I3").Value = "Ja"
$logs.Range("J3").Value = "Ja"

# Re-fit the row height: writing a multi-line string auto-grows the row,
# but the source row (2) uses the sheet's default height, so bring row 3
# back in line with it instead of leaving a stray explicit row height.
$logs.Rows.Item(3).AutoFit()

# --- Extend conditional formatting ranges to include row 3 -----------
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $srcCell = $col + "2"
    $newRange = $col + "2:" + $col + "3"
    $conditions = $logs.Range($srcCell).FormatConditions
    for ($i = 1; $i -le $conditions.Count; $i++) {
        $conditions.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# --- Bump the Dashboard "Overig" tally --------------------------------
$dash.Range("B2").Value = 2
